$d = $word.ActiveDocument

$pairs = @(
    @("50×63=", "33×95="),
    @("92×92=", "98×18="),
    @("75×44=", "20×27="),
    @("54×73=", "38×39="),
    @("99×33=", "14×93="),
    @("29×81=", "52×74="),
    @("21×34=", "48×50="),
    @("84×79=", "13×95="),
    @("39×65=", "40×96="),
    @("67×69=", "25×23="),
    @("62×27=", "34×68="),
    @("99×57=", "53×18="),
    @("24×20=", "37×23="),
    @("57×28=", "77×80="),
    @("63×26=", "48×85="),
    @("39×28=", "29×24="),
    @("14×20=", "96×30="),
    @("47×60=", "32×22="),
    @("54×26=", "33×51="),
    @("26×20=", "63×41="),
    @("47×52=", "89×24="),
    @("49×75=", "35×97="),
    @("99×11=", "22×42="),
    @("66×26=", "76×44="),
    @("96×39=", "54×37=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
